$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (changed) date column C for rows 2-5 from 45224 to 45233
$ws.Range("C2:C5").Value = 45233
